$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: columns H and I swap meaning.
# H1 was "Result" -> becomes "SamplePortion"
# I1 was "SamplePortion" -> becomes "Result"
$ws.Range("H1").ClearContents()
$ws.Range("I1").ClearContents()
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# Row 2: the float type marker now documents its unit (mg)
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# Row 3 (new): French descriptions for every column
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
